# 2017-07-19: 1. add dispatcher to router CRUD to different bussiness logic code.
#             2. add new part METHOD
#
# - Rename Sheet1 -> "part"
# - Add a "method" column (K) describing how each CRUD route is matched,
#   and rename the old "singleFieldValue" header (J1) to "singleField"
#   now that the "method" concept lives in its own column.
# - Add a review comment on the new K6 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- rename the first sheet ---
$ws.Name = "part"

# --- new "method" column header, then rename the old singleFieldValue header ---
$ws.Range("K1").Value = "method"
$ws.Range("J1").Value = "singleField"
$ws.Range("J1").Font.Name = "宋体"

# --- per-row method codes (K2:K5), matching the existing create/search/update/delete rows ---
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 1
$ws.Range("K5").Value = 3

# --- new row-6 entry documenting the 4th matching method ---
$ws.Range("K6").Value = "match:4"
$ws.Range("K6").Font.Name = "Arial"

$excel.UserName = "ZHANG Wei AG"
$ws.Range("K6").AddComment("ZHANG Wei AG:" + [char]10 + "判断输入记录是否和db中存储的一直（例如登录）")

# --- update the view selection to reflect where the edit happened ---
$ws.Activate()
$ws.Range("M6").Select()
